$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# --- Cells that flip from a numeric value to a text marker ("0" / "***.*") ---
# Copy number-format + style from an existing marker cell of the right style (C15 -> "0", L15 -> "***.*"),
# then write the text with a leading apostrophe so it is stored as a literal string, matching the source cells.
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = "'0"
$ws.Range("L15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = "'***.*"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").Value = "'0"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("D26").Value = "'0"
$ws.Range("L15").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").Value = "'***.*"

# --- Cells that flip from a text marker to a real numeric value ---
# Copy number-format + style from an existing numeric cell of the right style (F15 -> count style, H15 -> percent style),
# then write the number.
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Value = 1
$ws.Range("H15").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").Value = -100
$ws.Range("F15").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = 2
$ws.Range("F15").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").Value = 1
$ws.Range("H15").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = 100
$ws.Range("H15").Copy() | Out-Null
$ws.Range("M20").PasteSpecial(-4122) | Out-Null
$ws.Range("M20").Value = 200
$ws.Range("H15").Copy() | Out-Null
$ws.Range("L30").PasteSpecial(-4122) | Out-Null
$ws.Range("L30").Value = 0

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -41.176470588235
$ws.Range("I16").Value = 19
$ws.Range("J16").Value = 22
$ws.Range("K16").Value = -13.636363636363
$ws.Range("L16").Value = 11.764705882352
$ws.Range("M16").Value = 137.5
$ws.Range("N16").Value = -86.029411764705
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = -25
$ws.Range("J17").Value = 12
$ws.Range("K17").Value = 8.333333333333
$ws.Range("L17").Value = -23.529411764705
$ws.Range("N17").Value = -43.478260869565
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -71.428571428571
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -5.882352941176
$ws.Range("I18").Value = 24
$ws.Range("J18").Value = 34
$ws.Range("K18").Value = -29.411764705882
$ws.Range("L18").Value = 20
$ws.Range("M18").Value = -25
$ws.Range("N18").Value = -79.130434782608
$ws.Range("C19").Value = 28
$ws.Range("D19").Value = 36
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 85
$ws.Range("G19").Value = 107
$ws.Range("H19").Value = -20.560747663551
$ws.Range("I19").Value = 150
$ws.Range("J19").Value = 170
$ws.Range("K19").Value = -11.764705882352
$ws.Range("L19").Value = 92.307692307692
$ws.Range("M19").Value = 2.739726027397
$ws.Range("N19").Value = -68.152866242038
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 2
$ws.Range("K20").Value = 200
$ws.Range("N20").Value = -95.454545454545
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = -32
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 152
$ws.Range("H21").Value = -20.394736842105
$ws.Range("I21").Value = 213
$ws.Range("J21").Value = 242
$ws.Range("K21").Value = -11.983471074380
$ws.Range("L21").Value = 54.347826086956
$ws.Range("M21").Value = 8.121827411167
$ws.Range("N21").Value = -75.740318906605
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 7
$ws.Range("H22").Value = -22.222222222222
$ws.Range("I22").Value = 12
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = -7.692307692307
$ws.Range("L22").Value = -7.692307692307
$ws.Range("M22").Value = -7.692307692307
$ws.Range("C24").Value = 65
$ws.Range("D24").Value = 74
$ws.Range("E24").Value = -12.162162162162
$ws.Range("F24").Value = 287
$ws.Range("G24").Value = 292
$ws.Range("H24").Value = -1.712328767123
$ws.Range("I24").Value = 507
$ws.Range("J24").Value = 488
$ws.Range("K24").Value = 3.893442622950
$ws.Range("L24").Value = 102.8
$ws.Range("M24").Value = 156.060606060606
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 700
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = -4.166666666666
$ws.Range("I25").Value = 42
$ws.Range("J25").Value = 39
$ws.Range("K25").Value = 7.692307692307
$ws.Range("L25").Value = 40
$ws.Range("M25").Value = 35.483870967741
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 0
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = 75
$ws.Range("L27").Value = 100
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -75
$ws.Range("J30").Value = 7
$ws.Range("K30").Value = -85.714285714285

$excel.CutCopyMode = 0

Write-Output "done"
